$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date formatting already used by the table's Date column (style
# carries numFmtId 14, the built-in short-date format) onto the new rows
# before writing values, so we don't introduce a brand-new custom number
# format into styles.xml.
$ws.Range("A50").Copy()
$ws.Range("A51:A53").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Lesson 17: Thunderstorms videos watched.
$ws.Range("A51").Value = 45042
$ws.Range("B51").Value = "Your First Solo "
$ws.Range("C51").Value = "Closer Look: Get the Big Picture"

$ws.Range("A52").Value = 45042
$ws.Range("B52").Value = "Your First Solo "
$ws.Range("C52").Value = "Intro to Glass Cockpit Systems"

$ws.Range("A53").Value = 45047
$ws.Range("B53").Value = "Your First Solo "
$ws.Range("C53").Value = "Airport Signs and Markings"

# Match the author's final selection in the workbook.
$ws.Range("J38:J40").Select() | Out-Null
